# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: per-play yardage logs (space separated number lists) -
# append Week 16's logged plays.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 3 4 3 1 -1 7 7 4 4 0 2 2 10 3 1 2 19 0 2 1 2 4"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 4 17 8 11 12 4 14 -5 6 4 4 56"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 2 2 1 4 2 2 4 3 4 4 2 9 4 -3 1 10 6 2 0 4 4 1 5 3 4 0 0 5 -1 2"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 4 4 5 -1 18 7 11 6 7 6 10 9 17 13 40 24 1 5 12"

# ---------------------------------------------------------------------
# OFF sheet: season totals, row2 = Home, row3 = Road
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Range("C2").Value = 191
$offWs.Range("F2").Value = 67
$offWs.Range("G2").Value = 69
$offWs.Range("H2").Value = 10
$offWs.Range("J2").Value = 27
$offWs.Range("L2").Value = 246
$offWs.Range("M2").Value = 142
$offWs.Range("O2").Value = 20
$offWs.Range("P2").Value = 11
$offWs.Range("Q2").Value = 495

$offWs.Range("C3").Value = 149
$offWs.Range("E3").Value = 36
$offWs.Range("F3").Value = 91
$offWs.Range("H3").Value = 26
$offWs.Range("I3").Value = 55
$offWs.Range("J3").Value = 44
$offWs.Range("N3").Value = 16

# ---------------------------------------------------------------------
# DEF sheet: season totals, row2 = Home, row3 = Road
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Range("C2").Value = 184
$defWs.Range("D2").Value = 8
$defWs.Range("F2").Value = 58
$defWs.Range("G2").Value = 43
$defWs.Range("H2").Value = 6
$defWs.Range("I2").Value = 5
$defWs.Range("J2").Value = 29
$defWs.Range("L2").Value = 280
$defWs.Range("M2").Value = 187
$defWs.Range("Q2").Value = 474

$defWs.Range("B3").Value = 11
$defWs.Range("C3").Value = 160
$defWs.Range("E3").Value = 31
$defWs.Range("F3").Value = 109
$defWs.Range("G3").Value = 29
$defWs.Range("I3").Value = 70
$defWs.Range("J3").Value = 61
$defWs.Range("N3").Value = 16

# ---------------------------------------------------------------------
# ST sheet: special teams - both the season totals (row2/row3) and the
# per-game logs (D/RA/RM lists in columns B and D).
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Range("B2").Value = 68
$stWs.Range("D2").Value = 75
$stWs.Range("J2").Value = 43
$stWs.Range("K2").Value = 39
$stWs.Range("B3").Value = 42

$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 67"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 23"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 26"

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 49 42 46 45 38 42 54"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " -1 9 0 0 0 0 10"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 10 0 0 0 0 17"

# ---------------------------------------------------------------------
# TURNS sheet: turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Range("B2").Value = 9
$turnsWs.Range("C2").Value = 10
$turnsWs.Range("E2").Value = 6

$turnsWs.Range("D3").Value = 3
$turnsWs.Range("E3").Value = 5

# ---------------------------------------------------------------------
# PEN sheet: penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Range("D2").Value = 3
$penWs.Range("D4").Value = 11
